$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New time-tracking entries for 2023-09-19 .. 2023-09-28 (rows 179-187).
# Columns: A=Date B=StartTime C=EndTime D=Duration E=Hashtag F=Descriptor
#          G=IsSoftwareProject H=IsReleaseDay I=Year(formula) J=Month(formula)

$rows = @(
    @{ r=179; A="2023-09-19"; B="17:30"; C="18:00"; D="0h 30m"; E="#python"; F="nwreadinglistmanager v1.5.0";  G="True"; H="False" },
    @{ r=180; A="2023-09-20"; B="17:00"; C="17:45"; D="0h 45m"; E="#python"; F="nwreadinglistmanager v1.5.0";  G="True"; H="False" },
    @{ r=181; A="2023-09-21"; B="17:00"; C="17:45"; D="0h 45m"; E="#python"; F="nwreadinglistmanager v1.5.0";  G="True"; H="False" },
    @{ r=182; A="2023-09-22"; B="17:30"; C="20:30"; D="3h 00m"; E="#python"; F="nwreadinglistmanager v1.5.0";  G="True"; H="False" },
    @{ r=183; A="2023-09-22"; B="21:45"; C="00:45"; D="3h 00m"; E="#python"; F="nwreadinglistmanager v1.5.0";  G="True"; H="False" },
    @{ r=184; A="2023-09-23"; B="10:15"; C="11:15"; D="1h 00m"; E="#python"; F="nwreadinglistmanager v1.5.0";  G="True"; H="True"  },
    @{ r=185; A="2023-09-23"; B="13:15"; C="15:15"; D="2h 00m"; E="#python"; F="nwreadinglistmanager v1.5.0";  G="True"; H="True"  },
    @{ r=186; A="2023-09-27"; B="18:00"; C="20:30"; D="2h 30m"; E="#python"; F="nwtimetrackingmanager v1.0.0"; G="True"; H="False" },
    @{ r=187; A="2023-09-28"; B="17:45"; C="18:15"; D="0h 30m"; E="#python"; F="nwtimetrackingmanager v1.0.0"; G="True"; H="False" }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = "'" + $row.G
    $ws.Cells.Item($r, 8).Value = "'" + $row.H
    $ws.Cells.Item($r, 9).Formula = "=YEAR(A" + $r + ")"
    $ws.Cells.Item($r, 10).Formula = "=MONTH(A" + $r + ")"
}

$ws.Range("E177").Select() | Out-Null
